$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.6853526830673218
$ws.Range("B1").Value = 0.8333009481430054
$ws.Range("C1").Value = 3.560703754425049
$ws.Range("D1").Value = 2.181986570358276
$ws.Range("E1").Value = 0.9111688733100891
